# Apply cell-value updates per the diff (crypto price/volume refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.361.11"
$ws.Range("E2").Value = "'  +4.08%  "

$ws.Range("D3").Value = "'1.723.11"
$ws.Range("E3").Value = "'  +3.83%  "

$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "'  -0.12%  "

$ws.Range("D5").Value = "'239.89"
$ws.Range("E5").Value = "'  +1.99%  "

$ws.Range("E7").Value = "'  -1.00%  "

$ws.Range("D8").Value = "'0.2630"
$ws.Range("E8").Value = "'  +0.98%  "

$ws.Range("D9").Value = "'0.06242"
$ws.Range("E9").Value = "'  +1.67%  "

$ws.Range("D10").Value = "'1.716.95"
$ws.Range("E10").Value = "'  +3.45%  "

$ws.Range("D11").Value = "'0.07074"
$ws.Range("E11").Value = "'  +0.09%  "

$ws.Range("D12").Value = "'15.36"
$ws.Range("E12").Value = "'  +4.88%  "

$ws.Range("D13").Value = "'0.5939"
$ws.Range("E13").Value = "'  +0.54%  "

$ws.Range("D14").Value = "'4.410"
$ws.Range("E14").Value = "'  +0.64%  "

$ws.Range("D15").Value = "'76.36"
$ws.Range("E15").Value = "'  +2.86%  "

$ws.Range("D16").Value = "'0.9998"
$ws.Range("E16").Value = "'  -0.11%  "

$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "'  -0.08%  "

$ws.Range("D18").Value = "'26.351.50"
$ws.Range("E18").Value = "'  +4.10%  "

$ws.Range("E19").Value = "'  +1.16%  "

$ws.Range("E20").Value = "'  +1.84%  "

$ws.Range("D21").Value = "'1.937.24"
$ws.Range("E21").Value = "'  +3.43%  "

$ws.Range("E22").Value = "'  +2.76%  "

$ws.Range("D23").Value = "'8.785"
$ws.Range("E23").Value = "'  +1.82%  "

$ws.Range("D24").Value = "'5.342"
$ws.Range("E24").Value = "'  +0.41%  "

$ws.Range("D25").Value = "'134.86"
$ws.Range("E25").Value = "'  +0.97%  "

$ws.Range("D26").Value = "'15.23"
$ws.Range("E26").Value = "'  +1.40%  "

$ws.Range("E27").Value = "'  +0.72%  "

$ws.Range("D28").Value = "'1.767"
$ws.Range("E28").Value = "'  +5.14%  "

$ws.Range("D29").Value = "'106.69"
$ws.Range("E29").Value = "'  +2.40%  "

$ws.Range("D30").Value = "'4.023"
$ws.Range("E30").Value = "'  +1.92%  "

$ws.Range("D31").Value = "'3.701"
$ws.Range("E31").Value = "'  +2.22%  "

$ws.Range("D32").Value = "'0.07759"
$ws.Range("E32").Value = "'  +1.65%  "

$ws.Range("D33").Value = "'0.04466"
$ws.Range("E33").Value = "'  +2.65%  "

$ws.Range("D34").Value = "'2.609"
$ws.Range("E34").Value = "'  +0.12%  "

$ws.Range("D35").Value = "'0.9774"
$ws.Range("E35").Value = "'  +3.47%  "

$ws.Range("D36").Value = "'0.6214"
$ws.Range("E36").Value = "'  +1.27%  "

$ws.Range("D37").Value = "'115.97"
$ws.Range("E37").Value = "'  +18.87%  "

$ws.Range("E38").Value = "'  +7.93%  "

$ws.Range("D39").Value = "'2.412"
$ws.Range("E39").Value = "'  -7.36%  "

$ws.Range("B40").Value = "'RenderToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.912"
$ws.Range("E40").Value = "'  +4.66%  "

$ws.Range("B41").Value = "'PaxDollar"
$ws.Range("C41").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "'  +0.03%  "

$ws.Range("E42").Value = "'  -1.83%  "

$ws.Range("D43").Value = "'5.371"
$ws.Range("E43").Value = "'  +15.87%  "

$ws.Range("D44").Value = "'0.3830"
$ws.Range("E44").Value = "'  +1.96%  "

$ws.Range("D45").Value = "'0.1161"
$ws.Range("E45").Value = "'  +4.28%  "

$ws.Range("D46").Value = "'6.275"
$ws.Range("E46").Value = "'  +1.60%  "

$ws.Range("D47").Value = "'0.05290"
$ws.Range("E47").Value = "'  +0.77%  "

$ws.Range("D48").Value = "'30.60"
$ws.Range("E48").Value = "'  +4.23%  "

$ws.Range("D49").Value = "'7.659"
$ws.Range("E49").Value = "'  +4.94%  "

$ws.Range("D50").Value = "'0.3396"
$ws.Range("E50").Value = "'  +2.05%  "

$ws.Range("E51").Value = "'  +2.47%  "
